# "cambios de agosto, puntos fe de ratas e historico"
# - Update the "Nota" text in row 8 (column L)
# - Update the reporting period / notification dates in row 8 (columns B, C, J, K)
# - Shrink row 8's height now that the note text is shorter
# - Move the sheet's selection/view down to the new data entry point

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Nota" text (replaces the old legal citation with the new wording)
$ws.Range("L8").Value = "La Universidad Politécnica de Pachuca, no cuenta con las recomendaciones de casos especiales de organismos garantes de derechos humanos."

# Updated period dates (Excel date serials)
$ws.Range("B8").Value = 44652   # Fecha de inicio del periodo que se informa -> 2022-04-01
$ws.Range("C8").Value = 44742   # Fecha de término del periodo que se informa -> 2022-06-30
$ws.Range("J8").Value = 44753   # Fecha de validación -> 2022-07-11
$ws.Range("K8").Value = 44753   # Fecha de actualización -> 2022-07-11

# Row 8 is shorter now (90 -> 60 points) since the note text shrank
$ws.Rows("8:8").RowHeight = 60

# Update the active selection to reflect where editing left off
$ws.Range("L9").Select()
